$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting numeric-looking strings into real numbers, and without
# leaving a residual number-format style on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '34.333.48'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '1.803.56'
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").Value = '  +0.11%  '

Set-TextValue $ws.Range("D5") '227.89'
$ws.Range("E5").Value = '  +0.73%  '

Set-TextValue $ws.Range("D6") '0.577'
$ws.Range("E6").Value = '  +3.79%  '

$ws.Range("E7").Value = '  +0.11%  '

Set-TextValue $ws.Range("D8") '36.02'
$ws.Range("E8").Value = '  +9.45%  '

$ws.Range("E9").Value = '  +1.92%  '

Set-TextValue $ws.Range("D10") '0.0692'
$ws.Range("E10").Value = '  +0.33%  '

Set-TextValue $ws.Range("D11") '0.0967'
$ws.Range("E11").Value = '  +2.14%  '

$ws.Range("D12").Value = '2.063.50'
$ws.Range("E12").Value = '  +0.79%  '

Set-TextValue $ws.Range("D13") '11.65'
$ws.Range("E13").Value = '  +4.50%  '

$ws.Range("D14").Value = '1.799.50'
$ws.Range("E14").Value = '  +0.64%  '

Set-TextValue $ws.Range("D15") '0.644'
$ws.Range("E15").Value = '  +1.52%  '

$ws.Range("E16").Value = '  +4.56%  '

$ws.Range("D17").Value = '34.338.83'
$ws.Range("E17").Value = '  -0.04%  '

Set-TextValue $ws.Range("D18") '69.10'
$ws.Range("E18").Value = '  +0.86%  '

Set-TextValue $ws.Range("D19") '245.47'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").Value = '0.0₃0794'
$ws.Range("E20").Value = '  -0.27%  '

Set-TextValue $ws.Range("D21") '11.50'
$ws.Range("E21").Value = '  +2.47%  '

$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("E23").Value = '  +0.68%  '

Set-TextValue $ws.Range("D24") '172.81'
$ws.Range("E24").Value = '  +2.91%  '

$ws.Range("E25").Value = '  +3.13%  '

Set-TextValue $ws.Range("D26") '7.96'
$ws.Range("E26").Value = '  +8.63%  '

Set-TextValue $ws.Range("D27") '16.88'
$ws.Range("E27").Value = '  +1.90%  '

$ws.Range("E28").Value = '  +2.70%  '

$ws.Range("E29").Value = '  +0.08%  '

Set-TextValue $ws.Range("D30") '4.04'
$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("E31").Value = '  +1.00%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D32") '1.25'
$ws.Range("E32").Value = '  +1.46%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D33") '3.84'
$ws.Range("E33").Value = '  +0.93%  '

$ws.Range("E34").Value = '  +0.21%  '

$ws.Range("D35").Value = '1.393.85'
$ws.Range("E35").Value = '  -1.09%  '

Set-TextValue $ws.Range("D36") '0.673'
$ws.Range("E36").Value = '  -1.54%  '

Set-TextValue $ws.Range("D37") '2.46'
$ws.Range("E37").Value = '  -5.91%  '

$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("E39").Value = '  -0.32%  '

Set-TextValue $ws.Range("D41") '0.962'
$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("E42").Value = '  +0.87%  '

Set-TextValue $ws.Range("D43") '81.81'
$ws.Range("E43").Value = '  -2.89%  '

Set-TextValue $ws.Range("D44") '2.42'
$ws.Range("E44").Value = '  +0.38%  '

Set-TextValue $ws.Range("D45") '13.55'
$ws.Range("E45").Value = '  -2.29%  '

$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("E47").Value = '  -4.89%  '

$ws.Range("D48").Value = '1.964.06'
$ws.Range("E48").Value = '  +0.92%  '

Set-TextValue $ws.Range("D49") '104.97'
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("E51").Value = '  -0.16%  '
